$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (including the shared border/alignment style) from the
# existing data cell B4 down into the new B5 cell, then set its text.
$ws.Range("B4").Copy($ws.Range("B5"))
$ws.Range("B5").Value = "暂时搁置"

# Update the active selection to B5, matching the recorded view state
$ws.Range("B5").Select()
